$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 18.5
$ws.Range("B3").Value = 647500
$ws.Range("B12").Value = 1982142.857142857
$ws.Range("B29").Value = 2829642.857142857
$ws.Range("B31").Value = 2829642.857142857
